$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: copy formatting from row 20's relevant cells (A, B, D:F) ---
$ws.Range("A20").Copy()
$ws.Range("A21").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B20").Copy()
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("D20:F20").Copy()
$ws.Range("D21:F21").PasteSpecial(-4122)

# --- Row 22: same formatting source ---
$ws.Range("A20").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("B20").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("D20:F20").Copy()
$ws.Range("D22:F22").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Row 21 values: Construct Binary Tree from Preorder and Inorder Traversal ---
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "Construct Binary Tree from Preorder and Inorder Traversal"
$ws.Range("D21").Value = "Tree"
$ws.Range("E21").Value = "medium"
$ws.Range("F21").Value = "leetcode 105"

# --- Row 22 values: Closest Binary Search Tree Value II ---
$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "Closest Binary Search Tree Value II"
$ws.Range("D22").Value = "Tree"
$ws.Range("E22").Value = "hard"
$ws.Range("F22").Value = "leetcode 272"

$ws.Range("F28").Select()
